# "Begin fix time to string" - add new ticket rows (3-8) and correct the
# amount for the existing Kirov->Moskva row. The "Время" (time) column
# values are written as explicit strings (e.g. "1:7", "0:0", "3:2", "3:3")
# rather than letting Excel coerce them into time/date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 2 amount (64 -> 62)
$ws.Range("H2").Value = 62

# New data rows to append (index column A, id column B, start, end, price,
# time-as-text, ticket type, amount)
$rows = @(
    @(1, 2, "Калуга",     "Москва",           600,  "1:7", "Плацкарт", 100),
    @(2, 3, "Москва",     "Калуга",           700,  "0:0", "Плацкарт", 700),
    @(3, 4, "Москва",     "Хабаровск",        5000, "0:0", "Плацкарт", 100),
    @(4, 5, "Хабаровск",  "Москва",           4500, "3:2", "Плацкарт", 100),
    @(5, 6, "Москва",     "Владивосток",      3700, "3:3", "Плацкарт", 100),
    @(6, 7, "Москва",     "Санкт-Петербург",  3000, "0:0", "Плацкарт", 100)
)

$r = 3
foreach ($row in $rows) {
    # Copy the formatting/style of A2 (bordered, bold, centered index style)
    # down onto each new index cell in column A.
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
    $ws.Range("D" + $r).Value = $row[3]
    $ws.Range("E" + $r).Value = $row[4]
    $ws.Range("F" + $r).Value = $row[5]
    $ws.Range("G" + $r).Value = $row[6]
    $ws.Range("H" + $r).Value = $row[7]

    $r = $r + 1
}
